$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "div_po": rename two description values and add borders
# (this also causes "school"/"road" shared strings to be dropped and
# "community home"/"primary school" to be created, matching the diff)
# ------------------------------------------------------------------
$wsPo = $wb.Worksheets.Item("div_po")

$wsPo.Range("A1:B6").Borders.Color = -16777216
$wsPo.Range("A1:B6").Borders.LineStyle = 1

$wsPo.Range("B6").Value = "community home"
$wsPo.Range("B3").Value = "primary school"

$wsPo.Columns.Item(2).AutoFit()

# ------------------------------------------------------------------
# Sheet "RivisionHistory": append a new row to the revision table
# ------------------------------------------------------------------
$wsRev = $wb.Worksheets.Item("RivisionHistory")
$lo = $wsRev.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$wsRev.Range("A4:E4").Copy()
$wsRev.Range("A5:E5").PasteSpecial(-4122)

$wsRev.Range("A5").Value = 43354
$wsRev.Range("B5").Value = "0.3"
$wsRev.Range("C5").Value = "div_cw and div_po updated"
$wsRev.Range("D5").Value = "rishi"
$wsRev.Range("E5").Value = "rishi"

# ------------------------------------------------------------------
# Sheet "div_cw": update the step values and drop the trailing rows
# ------------------------------------------------------------------
$wsCw = $wb.Worksheets.Item("div_cw")

$wsCw.Range("D69").Value = "2015-04-30"
$wsCw.Range("D70").ClearContents()

$wsCw.Range("A71:E74").Delete(-4162)

# ------------------------------------------------------------------
# Restore the selection on each sheet so the saved view matches
# ------------------------------------------------------------------
$wsPo.Range("C12").Select()
$wsCw.Range("H67").Select()
$wsRev.Range("F6").Select()

Write-Output "edit applied"
